$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.799.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.759.03'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.49'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4253'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3620'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07573'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.54'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.093'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.004'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.65'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -6.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.054'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.265'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.775.97'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.11'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06379'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.04'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.905'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.833.59'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.20'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.115'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +6.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.28'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.27'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.970.21'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.44%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -8.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.91'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.115'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -7.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.686'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.578'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08889'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.48%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.63%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2106'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06026'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.85%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.982'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.86%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6344'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.180'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.891'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.401'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.40'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5876'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.702'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.983'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.83'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.179'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06833'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.09%  '
